# Auto-generated update of cryptos.xlsx crypto price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    if ($text -match '^[+-]?(\d+\.?\d*|\.\d+)([eE][+-]?\d+)?$') {
        # Value looks like a plain number to Excel's parser -- force
        # text interpretation (quote-prefix entry), then drop the
        # quote-prefix style so the cell keeps the workbook default style.
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range("D2") "67.348.76"
Set-TextValue $ws.Range("E2") "  -1.57%  "

Set-TextValue $ws.Range("D3") "3.749.69"
Set-TextValue $ws.Range("E3") "  +0.06%  "

Set-TextValue $ws.Range("E4") "  -0.07%  "

Set-TextValue $ws.Range("D5") "593.62"
Set-TextValue $ws.Range("E5") "  -0.28%  "

Set-TextValue $ws.Range("D6") "165.58"
Set-TextValue $ws.Range("E6") "  -0.82%  "

Set-TextValue $ws.Range("D7") "3.748.63"
Set-TextValue $ws.Range("E7") "  +0.04%  "

Set-TextValue $ws.Range("E8") "  +0.05%  "

Set-TextValue $ws.Range("D9") "0.518"
Set-TextValue $ws.Range("E9") "  -0.30%  "

Set-TextValue $ws.Range("E10") "  -0.50%  "

Set-TextValue $ws.Range("E11") "  -1.92%  "

Set-TextValue $ws.Range("E12") "  +0.10%  "

Set-TextValue $ws.Range("E13") "  -1.92%  "

Set-TextValue $ws.Range("D14") "36.04"
Set-TextValue $ws.Range("E14") "  +0.20%  "

Set-TextValue $ws.Range("D15") "4.379.49"
Set-TextValue $ws.Range("E15") "  +0.01%  "

Set-TextValue $ws.Range("D16") "3.751.66"
Set-TextValue $ws.Range("E16") "  -0.11%  "

Set-TextValue $ws.Range("D17") "18.37"
Set-TextValue $ws.Range("E17") "  +2.35%  "

Set-TextValue $ws.Range("D18") "67.325.09"
Set-TextValue $ws.Range("E18") "  -1.64%  "

Set-TextValue $ws.Range("E19") "  +0.11%  "

Set-TextValue $ws.Range("D20") "6.98"
Set-TextValue $ws.Range("E20") "  -0.26%  "

Set-TextValue $ws.Range("D21") "9.97"
Set-TextValue $ws.Range("E21") "  -7.21%  "

Set-TextValue $ws.Range("D22") "454.92"
Set-TextValue $ws.Range("E22") "  -2.20%  "

Set-TextValue $ws.Range("E23") "  -0.34%  "

Set-TextValue $ws.Range("D24") "0.0000154"
Set-TextValue $ws.Range("E24") "  +6.36%  "

Set-TextValue $ws.Range("D25") "83.14"
Set-TextValue $ws.Range("E25") "  -1.59%  "

Set-TextValue $ws.Range("D26") "2.14"
Set-TextValue $ws.Range("E26") "  -2.19%  "

Set-TextValue $ws.Range("D27") "11.86"
Set-TextValue $ws.Range("E27") "  -0.98%  "

Set-TextValue $ws.Range("D28") "10.12"
Set-TextValue $ws.Range("E28") "  +1.11%  "

Set-TextValue $ws.Range("E29") "  +0.06%  "

Set-TextValue $ws.Range("E30") "  -0.35%  "

Set-TextValue $ws.Range("D31") "7.25"
Set-TextValue $ws.Range("E31") "  -0.35%  "

Set-TextValue $ws.Range("D32") "29.57"
Set-TextValue $ws.Range("E32") "  -0.72%  "

Set-TextValue $ws.Range("D33") "2.18"
Set-TextValue $ws.Range("E33") "  +0.29%  "

Set-TextValue $ws.Range("D34") "9.16"
Set-TextValue $ws.Range("E34") "  -0.44%  "

Set-TextValue $ws.Range("E35") "  +0.01%  "

Set-TextValue $ws.Range("D36") "3.703.19"
Set-TextValue $ws.Range("E36") "  +0.01%  "

Set-TextValue $ws.Range("D37") "0.100"
Set-TextValue $ws.Range("E37") "  -0.59%  "

Set-TextValue $ws.Range("E38") "  -1.60%  "

Set-TextValue $ws.Range("E39") "  -1.06%  "

Set-TextValue $ws.Range("E40") "  -0.69%  "

Set-TextValue $ws.Range("E41") "  -1.23%  "

Set-TextValue $ws.Range("D42") "0.999"
Set-TextValue $ws.Range("E42") "  -0.11%  "

Set-TextValue $ws.Range("D44") "45.17"
Set-TextValue $ws.Range("E44") "  +3.18%  "

Set-TextValue $ws.Range("D45") "0.298"
Set-TextValue $ws.Range("E45") "  -1.89%  "

Set-TextValue $ws.Range("D46") "46.97"
Set-TextValue $ws.Range("E46") "  +2.28%  "

Set-TextValue $ws.Range("D47") "148.57"
Set-TextValue $ws.Range("E47") "  +1.58%  "

Set-TextValue $ws.Range("D48") "8.32"
Set-TextValue $ws.Range("E48") "  -2.65%  "

Set-TextValue $ws.Range("E49") "  -4.50%  "

Set-TextValue $ws.Range("D50") "389.31"
Set-TextValue $ws.Range("E50") "  +0.06%  "

Set-TextValue $ws.Range("D51") "26.09"
Set-TextValue $ws.Range("E51") "  +1.85%  "
